# Add a new column L capturing the latest LDLC price-check snapshot.
# L mirrors the structure of the existing per-timestamp columns: a bold,
# bordered, centered header in row 1 and either a copied price (rows with
# data) or a blank cell (rows not yet scraped) below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: new timestamp, formatted exactly like the other header cells.
$ws.Cells.Item(1, 12).Value = "2026-01-28 02:57:26"
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)

# Data rows 2-100: carry the last known price forward into the new column.
for ($r = 2; $r -le 100; $r++) {
    $ws.Cells.Item($r, 12).Value = $ws.Cells.Item($r, 11).Value2
}

# Rows 101-204: no price has been recorded yet for these products, so the
# new column stays blank there too (matching columns D-K on those rows).
for ($r = 101; $r -le 204; $r++) {
    $ws.Cells.Item($r, 12).Style = "Normal"
}
